# Atualizei dados bibi e add
# Update retention metrics on the active sheet to reflect the latest
# numbers for several cohort periods (num_customers / cohort_size /
# retention_rate columns C, D, E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: cohort 2021, period_index 5 -> num_customers 34 -> 35
$ws.Range("C22").Value = 35
$ws.Range("E22").Value = 0.01318764129615674

# Row 27: cohort 2022, period_index 4 -> num_customers 60 -> 61
$ws.Range("C27").Value = 61
$ws.Range("E27").Value = 0.02708703374777975

# Row 31: cohort 2023, period_index 3 -> num_customers 65 -> 67
$ws.Range("C31").Value = 67
$ws.Range("E31").Value = 0.02897923875432526

# Row 36: cohort 2024, period_index 1 -> num_customers 155 -> 156
$ws.Range("C36").Value = 156
$ws.Range("E36").Value = 0.08082901554404145

# Row 37: cohort 2025, period_index 0 -> num_customers/cohort_size 1038 -> 1052
$ws.Range("C37").Value = 1052
$ws.Range("D37").Value = 1052
